$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# 1) Update Status for case 17 (row 18, column I) from
#    "Verificar nova implementação." to "Mudança Implementada, Testar".
#    This also causes "Verificar nova implementação." to become an unused
#    shared string, which gets dropped from sharedStrings.xml automatically.
$ws.Range("I18").Value2 = "Mudança Implementada, Testar"

# 2) Apply an AutoFilter on column I (Status, the 9th column of A:K) so only
#    initiatives that are not yet fully implemented/tested remain visible.
#    This hides every row whose Status isn't one of the listed values.
$ws.Range("A1:K31").AutoFilter(9, @("Implementar batentes","Implementar calculo flexível também do FAP (atualmente está como obrigatório).","Pendente.","Testar erros de inserção de dados.","Verificar erros quando os dados não são numéricos, ou os parâmetros são incoerentes.","Verificar nova implementação."), 7)

# Row 18 now has a Status not included in the filter values above, so
# AutoFilter would hide it too; keep it visible explicitly.
$ws.Rows.Item(18).Hidden = $false

# 3) Update the view: scroll the frozen pane back to the top (A2) and move
#    the active selection to F20.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("F20").Select()

Write-Output "done"
